# Updated cryptos list values (Price / Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an untouched default-styled cell so every rewritten cell keeps its
# original (default) style instead of picking up a new "quote prefix" style.
$refStyle = $ws.Range("B2").Style

$updates = [ordered]@{
    "D2" = "26.256.85"
    "E2" = "  -0.06%  "
    "D3" = "1.592.55"
    "E4" = "  -0.03%  "
    "D5" = "213.09"
    "E7" = "  +0.01%  "
    "E8" = "  -0.47%  "
    "E9" = "  -0.36%  "
    "D10" = "18.98"
    "E10" = "  -1.86%  "
    "D12" = "1.816.96"
    "E12" = "  +0.23%  "
    "D13" = "1.604.66"
    "E13" = "  +1.05%  "
    "E14" = "  -0.92%  "
    "E15" = "  -2.37%  "
    "D16" = "63.83"
    "E16" = "  -0.92%  "
    "D17" = "26.254.71"
    "D19" = "215.71"
    "E19" = "  +0.87%  "
    "E20" = "  -1.27%  "
    "E21" = "  -0.03%  "
    "D22" = "4.30"
    "E22" = "  +0.17%  "
    "E23" = "  +0.53%  "
    "E24" = "  -2.34%  "
    "D25" = "145.24"
    "E25" = "  +0.10%  "
    "E27" = "  -1.27%  "
    "E28" = "  +0.77%  "
    "D29" = "15.12"
    "E29" = "  -0.46%  "
    "E30" = "  -1.28%  "
    "D31" = "1.15"
    "E31" = "  +0.15%  "
    "E32" = "  -0.37%  "
    "D33" = "1.419.90"
    "E33" = "  +5.85%  "
    "E34" = "  -0.37%  "
    "E35" = "  -0.84%  "
    "E36" = "  -1.24%  "
    "D37" = "0.571"
    "E37" = "  -4.13%  "
    "E38" = "  -0.66%  "
    "D39" = "0.826"
    "E39" = "  +1.20%  "
    "E40" = "  +0.08%  "
    "E41" = "  +0.03%  "
    "D42" = "0.936"
    "E42" = "  -9.19%  "
    "E43" = "  +0.73%  "
    "E44" = "  -0.06%  "
    "D45" = "1.729.16"
    "E45" = "  +0.30%  "
    "D46" = "60.92"
    "E46" = "  -1.49%  "
    "D47" = "86.72"
    "E47" = "  -1.53%  "
    "E48" = "  -0.96%  "
    "E49" = "  -0.31%  "
    "D50" = "0.0952"
    "E50" = "  -2.71%  "
    "D51" = "0.999"
    "E51" = "  +0.05%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces text entry so numeric-looking strings (e.g. '4.30',
    # '0.999') keep their exact text representation instead of being parsed as numbers.
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = $refStyle
}
